$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/formatting of column P (the last existing data column)
# onto the new column Q, covering the header row and all data rows.
$ws.Range("P1:P11").Copy()
$ws.Range("Q1:Q11").PasteSpecial(-4122)

# New header: "Дата вывода из эксплуатации" (Decommissioning date)
$ws.Range("Q1").Value = "Дата вывода из эксплуатации"

# Only the first data row gets a sample decommissioning date value.
$ws.Range("Q2").Value = "2022-01-25"

# Match the column width used for the new column in the target workbook
# (~24.17 characters, same rounding granularity as the other custom column widths).
$ws.Columns("Q").ColumnWidth = 23.3
